# Append the 2025-09-18 "Pick 3" / "Pick 4" results rows to the Results sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make the sheet view explicitly left-to-right (matches the authoring tool's output).
try { $excel.ActiveWindow.DisplayRightToLeft = $false } catch {}
try { $ws.DisplayRightToLeft = $false } catch {}

$newRows = @(
    @{ Row = 4; Date = "2025-09-18"; Game = "Pick 3"; Phase = "250918"; Result = "8-6-1";   InsertedAt = "2025-09-18T22:06:26.310+04:00" },
    @{ Row = 5; Date = "2025-09-18"; Game = "Pick 4"; Phase = "250918"; Result = "6-0-7-6"; InsertedAt = "2025-09-18T22:06:26.310+04:00" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    # Prefix the numeric-looking values with an apostrophe so Excel keeps them
    # as literal text instead of coercing them into dates / numbers.
    $ws.Cells.Item($row, 1).Value = "'" + $r.Date
    $ws.Cells.Item($row, 2).Value = $r.Game
    $ws.Cells.Item($row, 3).Value = "'" + $r.Phase
    $ws.Cells.Item($row, 4).Value = $r.Result
    $ws.Cells.Item($row, 5).Value = $r.InsertedAt
}

Write-Output "Added rows 4-5 to '$($ws.Name)'; used range is now $($ws.UsedRange.Address())"
